$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a weekly price series for "Palta" (avocado) sold at the
# "Feria Lagunitas de Puerto Montt" market. A new week's worth of data
# (two quality grades: Primera / Segunda) is being recorded at the top of
# the data block that starts at row 232, pushing the existing rows
# (232-245) down by two rows (to 234-247).

$ws.Rows("232:233").Insert()

# New row 232: Primera
$ws.Range("A232").Value = 4
$ws.Range("B232").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C232").Value = "Los Lagos"
$ws.Range("D232").Value = 44516
$ws.Range("E232").Value = 10
$ws.Range("F232").Value = "Fruta"
$ws.Range("G232").Value = 100106
$ws.Range("H232").Value = "Oleaginosos"
$ws.Range("I232").Value = 100106002
$ws.Range("J232").Value = "Palta"
$ws.Range("K232").Value = "Hass"
$ws.Range("L232").Value = "Primera"
$ws.Range("M232").Value = 400
$ws.Range("N232").Value = 4000
$ws.Range("O232").Value = 4100
$ws.Range("P232").Value = 4050
$ws.Range("Q232").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R232").Value = "Provincia de Quillota"
$ws.Range("S232").Value = 4050
$ws.Range("T232").Value = 1

# New row 233: Segunda
$ws.Range("A233").Value = 4
$ws.Range("B233").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C233").Value = "Los Lagos"
$ws.Range("D233").Value = 44516
$ws.Range("E233").Value = 10
$ws.Range("F233").Value = "Fruta"
$ws.Range("G233").Value = 100106
$ws.Range("H233").Value = "Oleaginosos"
$ws.Range("I233").Value = 100106002
$ws.Range("J233").Value = "Palta"
$ws.Range("K233").Value = "Hass"
$ws.Range("L233").Value = "Segunda"
$ws.Range("M233").Value = 200
$ws.Range("N233").Value = 3500
$ws.Range("O233").Value = 3500
$ws.Range("P233").Value = 3500
$ws.Range("Q233").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R233").Value = "Provincia de Quillota"
$ws.Range("S233").Value = 3500
$ws.Range("T233").Value = 1
